$d = $word.ActiveDocument

# --- Part 1: split "Oprimir los botones y que  cumplan cada uno " into
#     three runs: "Al ", "Oprimir ", "cualquiera de los botones se ejecute "
$r1 = $d.Content.Duplicate
$found1 = $r1.Find.Execute("Oprimir los botones y que  cumplan cada uno ")
if (-not $found1) {
    throw "Could not find target text for part 1"
}
$insertPos = $r1.Start
$r1.Text = ""

$pieceC = $d.Range($insertPos, $insertPos)
$pieceC.InsertBefore("cualquiera de los botones se ejecute ")

$pieceB = $d.Range($insertPos, $insertPos)
$pieceB.InsertBefore("Oprimir ")

$pieceA = $d.Range($insertPos, $insertPos)
$pieceA.InsertBefore("Al ")

# --- Part 2: split "su función específica la cual fue programada." into
#     "su función específica la cu" + "al fue programada." (bookmark stays
#     untouched right before this run)
$r2 = $d.Content.Duplicate
$found2 = $r2.Find.Execute("su función específica la cual fue programada.")
if (-not $found2) {
    throw "Could not find target text for part 2"
}
$splitOffset = $r2.Start + ("su función específica la cu".Length)
$tail = $d.Range($splitOffset, $r2.End)
$tail.Text = ""

$tailIns = $d.Range($splitOffset, $splitOffset)
$tailIns.InsertAfter("al fue programada.")
